$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I1 date label (stored as text)
$ws.Range("I1").Value = "14/03/2023"

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 889
$ws.Range("D2").Value = 1009
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 13
$ws.Range("H2").Value = 107
$ws.Range("I2").Value = 1008.8
$ws.Range("J2").Value = 0.01982553528945274
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 76
$ws.Range("D3").Value = 83
$ws.Range("E3").Value = 5
$ws.Range("I3").Value = 120
$ws.Range("J3").Value = -30.83333333333334
$ws.Range("C4").Value = 256
$ws.Range("D4").Value = 271
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 8
$ws.Range("I4").Value = 112
$ws.Range("J4").Value = 141.9642857142857
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 505
$ws.Range("D5").Value = 547
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 39
$ws.Range("I5").Value = 678
$ws.Range("J5").Value = -19.32153392330384
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 498
$ws.Range("D6").Value = 519
$ws.Range("E6").Value = 18
$ws.Range("G6").Value = 5
$ws.Range("I6").Value = 536
$ws.Range("J6").Value = -3.171641791044777
$ws.Range("C7").Value = 130
$ws.Range("D7").Value = 132
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("I7").Value = 165
$ws.Range("J7").Value = -20
$ws.Range("C8").Value = 98
$ws.Range("D8").Value = 131
$ws.Range("E8").Value = 33
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("I8").Value = 124
$ws.Range("J8").Value = 5.645161290322576
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 219
$ws.Range("D9").Value = 222
$ws.Range("E9").Value = 1
$ws.Range("G9").Value = 8
$ws.Range("I9").Value = 591
$ws.Range("J9").Value = -62.43654822335025
$ws.Range("C10").Value = 52
$ws.Range("D10").Value = 51
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("I10").Value = 86
$ws.Range("J10").Value = -40.69767441860465
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 4
$ws.Range("J11").Value = 33.33333333333333
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 1
$ws.Range("I12").Value = 68
$ws.Range("J12").Value = -70.58823529411764
